# Updated cryptos list with latest price and volume(1h) data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so numeric-looking values
# (e.g. "1.001", "0.9999") are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "23.763.22"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "1.655.17"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "303.61"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.3813"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("D8").Value = "0.3633"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "51.06"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "1.253"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").Value = "0.08208"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "22.69"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "6.546"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "7.465"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "0.00001240"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "1.652.84"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "97.69"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("D19").Value = "0.06999"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "6.802"
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("D21").Value = "17.79"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "12.83"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").Value = "23.765.00"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "2.530"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "3.050"
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").Value = "21.31"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "151.39"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "5.240"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "134.08"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").Value = "1.837.63"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").Value = "6.974"
$ws.Range("E32").Value = "  +4.35%  "
$ws.Range("D33").Value = "2.231"
$ws.Range("E33").Value = "  +2.81%  "
$ws.Range("D34").Value = "1.069"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").Value = "11.74"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("D36").Value = "0.02818"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").Value = "0.2523"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "6.141"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").Value = "0.08828"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").Value = "0.07137"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "13.19"
$ws.Range("E41").Value = "  +8.36%  "
$ws.Range("D42").Value = "0.7065"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "1.342"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "15.98"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "0.6538"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").Value = "2.331"
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("D47").Value = "0.9999"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "3.955"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "0.07959"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "127.47"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "1.191"
$ws.Range("E51").Value = "  -0.37%  "

# Restore the default (unstyled) cell style now that values are stored as text.
$priceRange.Style = "Normal"
